# Regenerate the "K" column (column G) values in the save_data sheet.
# This mirrors the author's commit: "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals" -- the recalculated K
# values replace the previous Strike# derived values for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$updates = @{
    2  = 1
    4  = 4
    5  = 3
    6  = 8
    7  = 7
    8  = 8
    9  = 6
    10 = 3
    11 = 7
    12 = 7
    13 = 10
    14 = 7
    15 = 10
    16 = 3
    17 = 1
    18 = 12
    19 = 6
    20 = 12
    21 = 8
    22 = 6
    23 = 10
    24 = 4
    25 = 3
    26 = 6
    27 = 10
    28 = 8
    29 = 9
    30 = 10
    31 = 5
    32 = 8
    33 = 5
    34 = 4
    35 = 4
    36 = 5
    37 = 6
    38 = 6
    39 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
